# Generate Report for Handback
# Updates the timestamp cells in the handback-status workbook to reflect
# a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
# This shared string is also referenced by de-de!H2 (same original text),
# so both cells must be updated together to keep sharing one string.
$wsOverview.Range("G2").Value = "2016-08-18 09:07:29"
$wsDeDe.Range("H2").Value = "2016-08-18 09:07:29"

# zh-cn sheet: "Correspond Handoff Datetime" (column H, row 2)
$wsZhCn.Range("H2").Value = "2016-08-18 09:07:03"

# zh-cn sheet: "Correspond Handback DateTime" (column K, row 2)
$wsZhCn.Range("K2").Value = "2016-08-18 09:07:44"

# de-de sheet: "Correspond Handback DateTime" (column K, row 2)
$wsDeDe.Range("K2").Value = "2016-08-18 09:07:52"
